# Applies the F/G column numeric updates described by the commit diff
# (sheet1=展览, sheet2=演出, sheet3=本地生活, sheet4=全部类型)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 509
$ws.Range("G5").Value = 63
$ws.Range("F6").Value = 957
$ws.Range("F7").Value = 196
$ws.Range("F9").Value = 1033
$ws.Range("F10").Value = 823
$ws.Range("F15").Value = 281
$ws.Range("F17").Value = 501
$ws.Range("F18").Value = 1332
$ws.Range("F19").Value = 120
$ws.Range("F20").Value = 892
$ws.Range("F21").Value = 1184
$ws.Range("F22").Value = 2869
$ws.Range("F23").Value = 1412
$ws.Range("F25").Value = 192
$ws.Range("F26").Value = 1272
$ws.Range("F28").Value = 1013
$ws.Range("F29").Value = 357
$ws.Range("F30").Value = 3087
$ws.Range("F31").Value = 602
$ws.Range("F32").Value = 535
$ws.Range("F33").Value = 1399

$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 78
$ws.Range("F15").Value = 2

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 739

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 739
$ws.Range("F5").Value = 78
$ws.Range("F7").Value = 509
$ws.Range("G7").Value = 63
$ws.Range("F12").Value = 957
$ws.Range("F13").Value = 196
$ws.Range("F16").Value = 1033
$ws.Range("F17").Value = 823
$ws.Range("F27").Value = 281
$ws.Range("F29").Value = 501
$ws.Range("F30").Value = 1332
$ws.Range("F31").Value = 120
$ws.Range("F32").Value = 892
$ws.Range("F33").Value = 1184
$ws.Range("F34").Value = 2869
$ws.Range("F35").Value = 1412
$ws.Range("F37").Value = 192
$ws.Range("F38").Value = 1272
$ws.Range("F42").Value = 1013
$ws.Range("F43").Value = 357
$ws.Range("F44").Value = 3087
$ws.Range("F45").Value = 602
$ws.Range("F46").Value = 535
$ws.Range("F47").Value = 1399
$ws.Range("F48").Value = 2
